# Modify seed value for Max Growth Rate Constraint
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Max Growth" seed formulas to reference the prior-year columns
# (H/C) instead of the final-year columns (J/E).
$ws.Range("L6").Formula  = "=-H46/1000"
$ws.Range("L7").Formula  = "=-H46/1000"

$ws.Range("L9").Formula  = "=-C45/1000"
$ws.Range("L10").Formula = "=-C45/1000"
$ws.Range("L11").Formula = "=-C45/1000"

$ws.Range("L12").Formula = "=-C44/1000"
$ws.Range("L13").Formula = "=-C44/1000"
$ws.Range("L14").Formula = "=-C44/1000"

$ws.Range("L15").Formula = "=-C43/1000"
$ws.Range("L16").Formula = "=-C43/1000"
$ws.Range("L17").Formula = "=-C43/1000"

# Restore the view state (scroll position + active selection) captured in
# the saved workbook.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K21").Select()
